# Update the "Förändrad" (modified) date in column C for all data rows
# from 2023-09-13 (serial 45182) to 2023-09-15 (serial 45184).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C468").Value = 45184
